# Commit: "Fruta / hortaliza, semanal"
# A new daily price record (row) needs to be inserted into the weekly
# "Ají" (Hortaliza) sheet at row 110, pushing the existing rows 110..178
# down to 111..179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 110 (shifts rows 110:178 down to 111:179)
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new record's data
$ws.Range("A110").Value = 5
$ws.Range("B110").Value = "Macroferia Regional de Talca"
$ws.Range("C110").Value = "Maule"
$ws.Range("D110").Value = 44596
$ws.Range("E110").Value = 7
$ws.Range("F110").Value = 100112021
$ws.Range("G110").Value = "Ají"
$ws.Range("H110").Value = "Americana (o)"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 150
$ws.Range("K110").Value = 8000
$ws.Range("L110").Value = 8000
$ws.Range("M110").Value = 8000
$ws.Range("N110").Value = "`$/caja 14 kilos"
$ws.Range("O110").Value = "Región del Maule"
$ws.Range("P110").Value = 571
$ws.Range("Q110").Value = 14
$ws.Range("R110").Value = "Hortaliza"
